$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column F (dSF) values per re-pulled data / mean calculation fix
$ws.Range("F6").Value = -7
$ws.Range("F7").Value = -1
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = -5
